$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 30313.074
$ws.Range("I98").Value = 929.2222
$ws.Range("J98").Value = 89080.78
$ws.Range("K98").Value = 929.2222
$ws.Range("L98").Value = 89080.78
$ws.Range("M98").Value = 568.7778
$ws.Range("N98").Value = -92076.78
# Row 99: Rumor Has It / Commanding Craftsman's Tea
$ws.Range("H99").Value = 1958.4286
$ws.Range("I99").Value = 1958.4286
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5875.2858
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -4377.2858
# Row 103: Let Loose the Juice / Persimmon Tannin
$ws.Range("H103").Value = 947.3333
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 947.3333
$ws.Range("K103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("M103").Value = 2841.9999
$ws.Range("N103").Value = -4013.9999
# Row 109: A Time for Peace / Smilodonskin Codex
$ws.Range("H109").Value = 34314
$ws.Range("J109").Value = 34314
$ws.Range("L109").Value = 34314
$ws.Range("N109").Value = -37088
# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 3289.8667
$ws.Range("I113").Value = 3310.5715
$ws.Range("K113").Value = 3310.5715
$ws.Range("M113").Value = -56.57150000000001
# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 30313.074
$ws.Range("I122").Value = 929.2222
$ws.Range("J122").Value = 89080.78
$ws.Range("K122").Value = 2787.6666
$ws.Range("L122").Value = 267242.34
$ws.Range("M122").Value = -337.6666
$ws.Range("N122").Value = -272142.34
# Row 127: Liquid Competence / Competent Craftsman's Draught
$ws.Range("H127").Value = 1091.1875
$ws.Range("I127").Value = 574.1667
$ws.Range("J127").Value = 1401.4
$ws.Range("K127").Value = 1722.5001
$ws.Range("L127").Value = 4204.200000000001
$ws.Range("M127").Value = 3237.4999
$ws.Range("N127").Value = -14124.2
# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 1376.3158
$ws.Range("I129").Value = 2165.6667
$ws.Range("J129").Value = 1228.3125
$ws.Range("K129").Value = 6497.000100000001
$ws.Range("L129").Value = 3684.9375
$ws.Range("M129").Value = -1497.000100000001
$ws.Range("N129").Value = -13684.9375
# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 17242552
$ws.Range("I135").Value = 1033.25
$ws.Range("J135").Value = 100001840
$ws.Range("K135").Value = 9299.25
$ws.Range("L135").Value = 900016560
$ws.Range("M135").Value = -6764.25
$ws.Range("N135").Value = -900021630
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 3170.106
$ws.Range("I137").Value = 898.55554
$ws.Range("J137").Value = 4021.9375
$ws.Range("K137").Value = 2695.66662
$ws.Range("L137").Value = 12065.8125
$ws.Range("M137").Value = -145.66662
$ws.Range("N137").Value = -17165.8125
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2032.0405
$ws.Range("I138").Value = 2209.2
$ws.Range("J138").Value = 1966.4259
$ws.Range("K138").Value = 6627.599999999999
$ws.Range("L138").Value = 5899.2777
$ws.Range("M138").Value = -1487.599999999999
$ws.Range("N138").Value = -16179.2777

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 1176.1923
$ws.Range("I2").Value = 1096.7727
$ws.Range("K2").Value = 1096.7727
$ws.Range("M2").Value = -983.7727
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 3005.2273
$ws.Range("I61").Value = 1604
$ws.Range("K61").Value = 1604
$ws.Range("M61").Value = -1392
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1524
$ws.Range("J74").Value = 2367.2727
$ws.Range("L74").Value = 2367.2727
$ws.Range("N74").Value = -4115.2727
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1524
$ws.Range("J77").Value = 2367.2727
$ws.Range("L77").Value = 11836.3635
$ws.Range("N77").Value = -20572.3635
# Row 97: Ore for Me / High Steel Ingot
$ws.Range("H97").Value = 977.3333
$ws.Range("I97").Value = 829.9
$ws.Range("J97").Value = 1272.2
$ws.Range("K97").Value = 829.9
$ws.Range("L97").Value = 1272.2
$ws.Range("M97").Value = -333.9
$ws.Range("N97").Value = -2264.2
# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 1176.1923
$ws.Range("I116").Value = 1096.7727
$ws.Range("K116").Value = 1096.7727
$ws.Range("M116").Value = 1197.2273
# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2090.6047
$ws.Range("I122").Value = 2244.2964
$ws.Range("J122").Value = 1831.25
$ws.Range("K122").Value = 6732.889200000001
$ws.Range("L122").Value = 5493.75
$ws.Range("M122").Value = -4282.889200000001
$ws.Range("N122").Value = -10393.75
# Row 124: Ace of Gloves / High Durium Gauntlets of Fending
$ws.Range("H124").Value = 30666.666
$ws.Range("J124").Value = 30666.666
$ws.Range("L124").Value = 30666.666
$ws.Range("N124").Value = -40486.666
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 13515217
$ws.Range("I132").Value = 20834594
$ws.Range("J132").Value = 2519.923
$ws.Range("K132").Value = 62503782
$ws.Range("L132").Value = 7559.768999999999
$ws.Range("M132").Value = -62501252
$ws.Range("N132").Value = -12619.769
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 3005.2273
$ws.Range("I136").Value = 1604
$ws.Range("K136").Value = 4812
$ws.Range("M136").Value = -2262
# Row 139: Backing up My Words / Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 37814.125
$ws.Range("J139").Value = 37814.125
$ws.Range("L139").Value = 37814.125
$ws.Range("N139").Value = -48094.125

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 1176.1923
$ws.Range("I3").Value = 1096.7727
$ws.Range("K3").Value = 1096.7727
$ws.Range("M3").Value = -982.7727
# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 1626.1852
$ws.Range("I99").Value = 1652.4783
$ws.Range("J99").Value = 1475
$ws.Range("K99").Value = 1652.4783
$ws.Range("L99").Value = 1475
$ws.Range("M99").Value = -154.4783
$ws.Range("N99").Value = -4471
# Row 137: Dagger Swagger / Cobalt Tungsten Khukuri
$ws.Range("H137").Value = 39756.125
$ws.Range("J137").Value = 39756.125
$ws.Range("L137").Value = 39756.125
$ws.Range("N137").Value = -49956.125

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 250690.86
$ws.Range("I31").Value = 2411.0667
$ws.Range("K31").Value = 2411.0667
$ws.Range("M31").Value = -2116.0667
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 250690.86
$ws.Range("I34").Value = 2411.0667
$ws.Range("K34").Value = 2411.0667
$ws.Range("M34").Value = -2209.0667

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water / Boiled Egg
$ws.Range("H4").Value = 1128.5
$ws.Range("J4").Value = 2157
$ws.Range("L4").Value = 6471
$ws.Range("N4").Value = -6695
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 2308.825
$ws.Range("I131").Value = 14864.143
$ws.Range("J131").Value = 1104.8904
$ws.Range("K131").Value = 44592.429
$ws.Range("L131").Value = 3314.6712
$ws.Range("M131").Value = -39552.429
$ws.Range("N131").Value = -13394.6712

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 4246.517
$ws.Range("I132").Value = 1590.4667
$ws.Range("J132").Value = 7092.2856
$ws.Range("K132").Value = 4771.4001
$ws.Range("L132").Value = 21276.8568
$ws.Range("M132").Value = -2241.4001
$ws.Range("N132").Value = -26336.8568
# Row 135: Fan of the Foreign / Ruthenium Folding Fans
$ws.Range("H135").Value = 42599.8
$ws.Range("J135").Value = 42599.8
$ws.Range("L135").Value = 42599.8
$ws.Range("N135").Value = -52739.8

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 2536.2144
$ws.Range("I7").Value = 2188.2666
$ws.Range("J7").Value = 2937.6924
$ws.Range("K7").Value = 2188.2666
$ws.Range("L7").Value = 2937.6924
$ws.Range("M7").Value = -2076.2666
$ws.Range("N7").Value = -3161.6924
# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 2192.1875
$ws.Range("I40").Value = 2133.75
$ws.Range("J40").Value = 2367.5
$ws.Range("K40").Value = 2133.75
$ws.Range("L40").Value = 2367.5
$ws.Range("M40").Value = -1997.75
$ws.Range("N40").Value = -2639.5
# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 2945.9167
$ws.Range("I61").Value = 2536.1
$ws.Range("J61").Value = 4995
$ws.Range("K61").Value = 2536.1
$ws.Range("L61").Value = 4995
$ws.Range("M61").Value = -2334.1
$ws.Range("N61").Value = -5399
# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 5953301.5
$ws.Range("J82").Value = 11905449
$ws.Range("L82").Value = 11905449
$ws.Range("N82").Value = -11906171
# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 5953301.5
$ws.Range("J85").Value = 11905449
$ws.Range("L85").Value = 11905449
$ws.Range("N85").Value = -11907945
# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 1680.7142
$ws.Range("I93").Value = 1299.1
$ws.Range("K93").Value = 1299.1
$ws.Range("M93").Value = -51.09999999999991
# Row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 1907.8334
$ws.Range("I100").Value = 1940.1875
$ws.Range("J100").Value = 1649
$ws.Range("K100").Value = 1940.1875
$ws.Range("L100").Value = 1649
$ws.Range("M100").Value = -1399.1875
$ws.Range("N100").Value = -2731
# Row 111: Glove Me Tender / Gliderskin Gloves of Striking
$ws.Range("H111").Value = 43939.6
$ws.Range("J111").Value = 43939.6
$ws.Range("L111").Value = 43939.6
$ws.Range("N111").Value = -52119.6
# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 2945.9167
$ws.Range("I113").Value = 2536.1
$ws.Range("J113").Value = 4995
$ws.Range("K113").Value = 2536.1
$ws.Range("L113").Value = 4995
$ws.Range("M113").Value = -366.0999999999999
$ws.Range("N113").Value = -9335
# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 2283.5833
$ws.Range("I122").Value = 2254.889
$ws.Range("K122").Value = 6764.667
$ws.Range("M122").Value = -4314.667
# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 2536.2144
$ws.Range("I126").Value = 2188.2666
$ws.Range("J126").Value = 2937.6924
$ws.Range("K126").Value = 6564.7998
$ws.Range("L126").Value = 8813.0772
$ws.Range("M126").Value = -4094.7998
$ws.Range("N126").Value = -13753.0772
# Row 127: Loyal Turncoat / Saigaskin Coat of Fending
$ws.Range("H127").Value = 48276.332
$ws.Range("J127").Value = 48276.332
$ws.Range("L127").Value = 48276.332
$ws.Range("N127").Value = -58196.332
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 2424.1765
$ws.Range("I136").Value = 1961.32
$ws.Range("K136").Value = 5883.96
$ws.Range("M136").Value = -3333.96

$ws = $wb.Worksheets.Item("WVR")
# Row 119: A Job Well Done / Dwarven Cotton Gaskins of Fending
$ws.Range("H119").Value = 48698
$ws.Range("J119").Value = 48698
$ws.Range("L119").Value = 48698
$ws.Range("N119").Value = -58374
# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 4202508.5
$ws.Range("I126").Value = 4902676.5
$ws.Range("K126").Value = 14708029.5
$ws.Range("M126").Value = -14705559.5
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 20412.906
$ws.Range("I136").Value = 48184.953
$ws.Range("J136").Value = 2187.5
$ws.Range("K136").Value = 144554.859
$ws.Range("L136").Value = 6562.5
$ws.Range("M136").Value = -142004.859
$ws.Range("N136").Value = -11662.5
